$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.31038849023281
$ws.Range("C2").Value = 1.104468044484804
$ws.Range("D2").Value = 0.07833833886036246
$ws.Range("E2").Value = 0.4103847033970851
$ws.Range("G2").Value = 0.002456054195868198
$ws.Range("I2").Value = 1.380927956172741
$ws.Range("N2").Value = 1.558116747980961
$ws.Range("B3").Value = 2.087521168377748
$ws.Range("C3").Value = 0.9870370563596111
$ws.Range("D3").Value = 0.0710239932534904
$ws.Range("E3").Value = 0.3573130995708453
$ws.Range("G3").Value = 0.002464655211288213
$ws.Range("I3").Value = 1.326349477690698
$ws.Range("N3").Value = 1.550947099941951
$ws.Range("B4").Value = 1.952249894350871
$ws.Range("C4").Value = 0.9156467929354903
$ws.Range("D4").Value = 0.06658642826411665
$ws.Range("E4").Value = 0.3249288306432874
$ws.Range("G4").Value = 0.002470192079407378
$ws.Range("I4").Value = 1.29379966121499
$ws.Range("N4").Value = 1.547110235807608
$ws.Range("B5").Value = 1.897507106503781
$ws.Range("C5").Value = 0.8867257192262059
$ws.Range("D5").Value = 0.06479107234842729
$ws.Range("E5").Value = 0.3117779595192758
$ws.Range("G5").Value = 0.002472513047666716
$ws.Range("I5").Value = 1.2807719518146
$ws.Range("N5").Value = 1.545686433881627
$ws.Range("B6").Value = 1.88843975319935
$ws.Range("C6").Value = 0.881933490588608
$ws.Range("D6").Value = 0.06449372609381498
$ws.Range("E6").Value = 0.3095969064192872
$ws.Range("G6").Value = 0.002472902357025461
$ws.Range("I6").Value = 1.27862285105715
$ws.Range("N6").Value = 1.545458383697593
$ws.Range("B7").Value = 1.951510087463589
$ws.Range("C7").Value = 0.9152560714185256
$ws.Range("D7").Value = 0.06656216347776933
$ws.Range("E7").Value = 0.3247512934795083
$ws.Range("G7").Value = 0.002470223118560934
$ws.Range("I7").Value = 1.293623013751983
$ws.Range("N7").Value = 1.547090471216407
$ws.Range("B8").Value = 2.233209680455218
$ws.Range("C8").Value = 1.063824862137835
$ws.Range("D8").Value = 0.0758049732036028
$ws.Range("E8").Value = 0.39204077999932
$ws.Range("G8").Value = 0.002458966928937606
$ws.Range("I8").Value = 1.361906543676056
$ws.Range("N8").Value = 1.555525856803314
$ws.Range("B9").Value = 2.798728807883776
$ws.Range("C9").Value = 1.361215418627467
$ws.Range("D9").Value = 0.09437637035010482
$ws.Range("E9").Value = 0.5258233402816757
$ws.Range("G9").Value = 0.00243890838677278
$ws.Range("I9").Value = 1.503686168501176
$ws.Range("N9").Value = 1.576664718395079
$ws.Range("B10").Value = 3.223164467308266
$ws.Range("C10").Value = 1.583968751356451
$ws.Range("D10").Value = 0.1083252433422359
$ws.Range("E10").Value = 0.6255557338232904
$ws.Range("G10").Value = 0.002425378575158125
$ws.Range("I10").Value = 1.613011165616413
$ws.Range("N10").Value = 1.595156823469125
$ws.Range("B11").Value = 3.418411374456696
$ws.Range("C11").Value = 1.68635960220297
$ws.Range("D11").Value = 0.1147443807944342
$ws.Range("E11").Value = 0.6713136812428218
$ws.Range("G11").Value = 0.002419481058210324
$ws.Range("I11").Value = 1.663946860011706
$ws.Range("N11").Value = 1.604247826120911
$ws.Range("B12").Value = 3.492675701803364
$ws.Range("C12").Value = 1.725295184530637
$ws.Range("D12").Value = 0.117186333249748
$ws.Range("E12").Value = 0.6887030751585286
$ws.Range("G12").Value = 0.00241728445344004
$ws.Range("I12").Value = 1.683414478957303
$ws.Range("N12").Value = 1.607790914569534
$ws.Range("B13").Value = 3.476666666357232
$ws.Range("C13").Value = 1.716902324278067
$ws.Range("D13").Value = 0.1166599094114531
$ws.Range("E13").Value = 0.6849551114119379
$ws.Range("G13").Value = 0.002417755906711953
$ws.Range("I13").Value = 1.679213697883711
$ws.Range("N13").Value = 1.607023329258368
$ws.Range("B14").Value = 3.424514463187052
$ws.Range("C14").Value = 1.689559542676079
$ws.Range("D14").Value = 0.1149450546550099
$ws.Range("E14").Value = 0.6727430431353696
$ws.Range("G14").Value = 0.002419299609686816
$ws.Range("I14").Value = 1.665544840700036
$ws.Range("N14").Value = 1.604537286619717
$ws.Range("B15").Value = 3.392613024726302
$ws.Range("C15").Value = 1.672832751462238
$ws.Range("D15").Value = 0.1138961275341899
$ws.Range("E15").Value = 0.6652710316139263
$ws.Range("G15").Value = 0.002420249935939087
$ws.Range("I15").Value = 1.657195822940224
$ws.Range("N15").Value = 1.603027692225993
$ws.Range("B16").Value = 3.210449680910983
$ws.Range("C16").Value = 1.577299437156171
$ws.Range("D16").Value = 0.1079072688789893
$ws.Range("E16").Value = 0.6225736748866524
$ws.Range("G16").Value = 0.002425769139838607
$ws.Range("I16").Value = 1.60970714911312
$ws.Range("N16").Value = 1.594576611990988
$ws.Range("B17").Value = 3.099265592865663
$ws.Range("C17").Value = 1.518971211380631
$ws.Range("D17").Value = 0.1042525760041997
$ws.Range("E17").Value = 0.5964840407524861
$ws.Range("G17").Value = 0.002429220644852069
$ws.Range("I17").Value = 1.580887022411389
$ws.Range("N17").Value = 1.589568062606304
$ws.Range("B18").Value = 3.035518574501111
$ws.Range("C18").Value = 1.485521439008039
$ws.Range("D18").Value = 0.102157398008373
$ws.Range("E18").Value = 0.5815142417347232
$ws.Range("G18").Value = 0.002431230094333232
$ws.Range("I18").Value = 1.564423308164848
$ws.Range("N18").Value = 1.586750996199996
$ws.Range("B19").Value = 3.013969317733881
$ws.Range("C19").Value = 1.474212636091352
$ws.Range("D19").Value = 0.1014491757560592
$ws.Range("E19").Value = 0.5764517777492131
$ws.Range("G19").Value = 0.002431914632133688
$ws.Range("I19").Value = 1.558868170468799
$ws.Range("N19").Value = 1.585808042275175
$ws.Range("B20").Value = 3.111080186796698
$ws.Range("C20").Value = 1.525170024151976
$ws.Range("D20").Value = 0.1046409056673099
$ws.Range("E20").Value = 0.5992575311142616
$ws.Range("G20").Value = 0.00242885072050616
$ws.Range("I20").Value = 1.583943248003422
$ws.Range("N20").Value = 1.590094615780686
$ws.Range("B21").Value = 3.439823776994331
$ws.Range("C21").Value = 1.697586291460425
$ws.Range("D21").Value = 0.1154484414975627
$ws.Range("E21").Value = 0.6763282971273696
$ws.Range("G21").Value = 0.002418845194394919
$ws.Range("I21").Value = 1.66955479558186
$ws.Range("N21").Value = 1.605264745921545
$ws.Range("B22").Value = 3.656599934747646
$ws.Range("C22").Value = 1.811221712757174
$ws.Range("D22").Value = 0.1225771391655428
$ws.Range("E22").Value = 0.7270618861678031
$ws.Range("G22").Value = 0.002412519509766642
$ws.Range("I22").Value = 1.726555518486123
$ws.Range("N22").Value = 1.615766628526842
$ws.Range("B23").Value = 3.540720801741486
$ws.Range("C23").Value = 1.750481886616171
$ws.Range("D23").Value = 0.1187662492608439
$ws.Range("E23").Value = 0.6999491709989201
$ws.Range("G23").Value = 0.002415876223229637
$ws.Range("I23").Value = 1.696035062196728
$ws.Range("N23").Value = 1.610106856982469
$ws.Range("B24").Value = 3.105738267555921
$ws.Range("C24").Value = 1.522367280321191
$ws.Range("D24").Value = 0.1044653233243338
$ws.Range("E24").Value = 0.5980035447438468
$ws.Range("G24").Value = 0.002429017884971518
$ws.Range("I24").Value = 1.58256120071762
$ws.Range("N24").Value = 1.589856366899753
$ws.Range("B25").Value = 2.644240286449246
$ws.Range("C25").Value = 1.280057604078536
$ws.Range("D25").Value = 0.08930124219965307
$ws.Range("E25").Value = 0.489402572772633
$ws.Range("G25").Value = 0.002444121224015902
$ws.Range("I25").Value = 1.464451400502284
$ws.Range("N25").Value = 1.570437745741742
